$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    35.36355261064645,
    38.48840163139249,
    40.33928021310246,
    41.13250984825467,
    36.82726021617131,
    41.43643303789648,
    38.62649634902951,
    37.29151186388066,
    36.38762104345184,
    36.26878131643321,
    31.41493417656761,
    35.03963592425407,
    33.61008544654723,
    26.06827386074252,
    25.6504450875162,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923,
    23.07929688879923
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
